# Scoring work-in-progress: the "Parsing" column (B) values such as
# "singular nominative" / "plural ablative" lose their internal space
# (e.g. "singularnominative") so they can be matched/scored as single
# tokens. This is a pure text edit of existing cells; Excel itself takes
# care of restructuring the shared-string table (renamed strings move to
# the end of the table, the untouched "Ending" strings shift up to fill
# the gap) when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $text = $cell.Text
    if ($text -match " ") {
        $cell.Value = ($text -replace " ", "")
    }
}

# Leave the selection where the author's cursor ended up.
$ws.Range("G8:G9").Select()
